$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [hashtable]$values
    )
    foreach ($col in $values.Keys) {
        $addr = "$col$row"
        $val = $values[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 40 @{ H=3012.75; I=1024.5; J=5001; K=1024.5; L=5001; M=-849.5; N=-5351 }
Set-Row $ws 63 @{ H=0; J=0; L=0; N=$null }
Set-Row $ws 66 @{ H=0; J=0; L=0; N=$null }
Set-Row $ws 106 @{ H=35454.453; I=52142.715; J=6250; K=52142.715; L=6250; M=-51511.715; N=-7512 }
Set-Row $ws 112 @{ H=2024.5834; I=975; J=2120; K=2925; L=6360; M=-1817; N=-8576 }
Set-Row $ws 121 @{ H=587.0454999999999; I=0; J=587.0454999999999; K=0; L=1761.1365; M=$null; N=-5255.1365 }
Set-Row $ws 131 @{ H=4469.5293; I=483.06668; J=6130.5557; K=1449.20004; L=18391.6671; M=3590.79996; N=-28471.6671 }
Set-Row $ws 141 @{ H=1488.7222; I=1334.8235; K=4004.4705; M=1175.5295 }

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 32 @{ H=4374.2856; I=3789.77; J=8997.272000000001; K=3789.77; L=8997.272000000001; M=-3502.77; N=-9571.272000000001 }
Set-Row $ws 97 @{ H=1925; I=1287.1428; K=1287.1428; M=-791.1428000000001 }
Set-Row $ws 102 @{ H=1014.2857; I=997.8; J=1055.5; K=997.8; L=1055.5; M=624.2; N=-4299.5 }
Set-Row $ws 132 @{ H=149163.98; I=26065.438; J=360190.06; K=78196.314; L=1080570.18; M=-75666.314; N=-1085630.18 }

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 134 @{ H=18201610; I=1488.1459; J=143002450; K=4464.4377; L=429007350; M=-1929.4377; N=-429012420 }

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 60 @{ H=7163.2964; I=4750; J=8179.421; K=4750; L=8179.421; M=-4239; N=-9201.421 }

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 2 @{ H=1050.6875; I=1535.5; J=759.8; K=9213; L=4558.799999999999; M=-9100; N=-4784.799999999999 }
Set-Row $ws 38 @{ H=5; J=0; L=0; N=$null }

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 132 @{ H=49677.094; I=1514.7858; J=146001.72; K=4544.357400000001; L=438005.16; M=-2014.357400000001; N=-443065.16 }

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 16 @{ H=1425.5652; I=1275.3846; J=1620.8; K=1275.3846; L=1620.8; M=-1105.3846; N=-1960.8 }
Set-Row $ws 61 @{ H=2942.6924; I=3158.3333; J=2757.8572; K=3158.3333; L=2757.8572; M=-2956.3333; N=-3161.8572 }
Set-Row $ws 68 @{ H=1844.2858; I=1390; J=2980; K=1390; L=2980; M=-641; N=-4478 }
Set-Row $ws 71 @{ H=1844.2858; I=1390; J=2980; K=6950; L=14900; M=-3206; N=-22388 }
Set-Row $ws 93 @{ H=1489.871; I=1074; J=2148.3333; K=1074; L=2148.3333; M=174; N=-4644.3333 }
Set-Row $ws 113 @{ H=2942.6924; I=3158.3333; J=2757.8572; K=3158.3333; L=2757.8572; M=-988.3332999999998; N=-7097.8572 }
Set-Row $ws 122 @{ H=2808; I=2893.2354; J=2696.5386; K=8679.706200000001; L=8089.6158; M=-6229.706200000001; N=-12989.6158 }
Set-Row $ws 132 @{ H=181298.67; I=45427.89; J=478920.38; K=136283.67; L=1436761.14; M=-133753.67; N=-1441821.14 }

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

Set-Row $ws 132 @{ H=4624.6333; I=786.75; J=12300.4; K=2360.25; L=36901.2; M=169.75; N=-41961.2 }
